$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row-content swaps: two match rows exchange all of their match data
#    (columns F..V) while keeping their own Indice/pais/torneio/temporada/
#    data_partida (columns A..E) untouched. Simple two-row swaps.
# ---------------------------------------------------------------------------
function Swap-Rows([int]$r1, [int]$r2) {
    $range1 = $ws.Range("F$r1`:V$r1")
    $range2 = $ws.Range("F$r2`:V$r2")
    $v1 = $range1.Value()
    $v2 = $range2.Value()
    $range1.Value = $v2
    $range2.Value = $v1
}

Swap-Rows 47 48
Swap-Rows 56 57
Swap-Rows 66 67
Swap-Rows 214 215
Swap-Rows 216 218
Swap-Rows 219 220

# ---------------------------------------------------------------------------
# 2) Rows 83, 84, 85 rotate: 83 takes 84's old data, 84 takes 85's old data,
#    85 takes 83's old data.
# ---------------------------------------------------------------------------
$r83 = $ws.Range("F83:V83").Value()
$r84 = $ws.Range("F84:V84").Value()
$r85 = $ws.Range("F85:V85").Value()

$ws.Range("F83:V83").Value = $r84
$ws.Range("F84:V84").Value = $r85
$ws.Range("F85:V85").Value = $r83

# ---------------------------------------------------------------------------
# 3) Append three brand-new match rows (232, 233, 234) at the end of the
#    sheet, matching the formatting already used by the preceding rows
#    (bold/bordered/centered column A, date-formatted column E).
# ---------------------------------------------------------------------------
$ws.Range("A231:V231").Copy()
$ws.Range("A232:V234").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# NOTE: this runtime's PowerShell subset does not bind named (-paramName
# value) arguments reliably, so the helper below takes purely positional
# parameters.
function Set-MatchRow($row, $indice, $fecha, $home, $homeGoals, $away, $awayGoals,
                       $homeOpenOdds, $homeOpenDt, $homeCloseOdds, $homeCloseDt,
                       $drawOpenOdds, $drawOpenDt, $drawCloseOdds, $drawCloseDt,
                       $awayOpenOdds, $awayOpenDt, $awayCloseOdds, $awayCloseDt, $url) {

    $ws.Range("A$row").Value = $indice
    $ws.Range("B$row").Value = "brazil"
    $ws.Range("C$row").Value = "serie-a"

    # "2023" looks numeric, so a plain .Value assignment would silently turn
    # it into a number cell; the source file stores it as text. Force text
    # via a temporary "@" format, then paste the (already-correct, General)
    # format from C2 back over it so the cell style index is unaffected.
    $ws.Range("D$row").NumberFormat = "@"
    $ws.Range("D$row").Value = "2023"
    $ws.Range("C2").Copy()
    $ws.Range("D$row").PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    $ws.Range("E$row").Value = $fecha
    $ws.Range("F$row").Value = $home
    $ws.Range("G$row").Value = $homeGoals
    $ws.Range("H$row").Value = $away
    $ws.Range("I$row").Value = $awayGoals
    $ws.Range("J$row").Value = $homeOpenOdds
    $ws.Range("K$row").Value = $homeOpenDt
    $ws.Range("L$row").Value = $homeCloseOdds
    $ws.Range("M$row").Value = $homeCloseDt
    $ws.Range("N$row").Value = $drawOpenOdds
    $ws.Range("O$row").Value = $drawOpenDt
    $ws.Range("P$row").Value = $drawCloseOdds
    $ws.Range("Q$row").Value = $drawCloseDt
    $ws.Range("R$row").Value = $awayOpenOdds
    $ws.Range("S$row").Value = $awayOpenDt
    $ws.Range("T$row").Value = $awayCloseOdds
    $ws.Range("U$row").Value = $awayCloseDt
    $ws.Range("V$row").Value = $url
}

Set-MatchRow 232 231 45190 `
    "Goias" 0 "Flamengo RJ" 0 `
    3.39 "16/09/2023 01:42" 2.96 "20/09/2023 23:59" `
    3.45 "16/09/2023 01:42" 3.2 "20/09/2023 23:48" `
    2.2 "16/09/2023 01:42" 2.62 "20/09/2023 23:48" `
    "https://www.betexplorer.com/football/brazil/serie-a/goias-flamengo-rj/vkEpb1vk/"

Set-MatchRow 233 232 45190.10416666666 `
    "Fluminense" 1 "Cruzeiro" 0 `
    1.78 "16/09/2023 20:12" 1.83 "21/09/2023 02:20" `
    3.76 "16/09/2023 20:12" 3.59 "21/09/2023 02:20" `
    4.66 "16/09/2023 20:12" 4.9 "21/09/2023 02:22" `
    "https://www.betexplorer.com/football/brazil/serie-a/fluminense-cruzeiro/niZ55qnF/"

Set-MatchRow 234 233 45190.10416666666 `
    "Sao Paulo" 1 "Fortaleza" 2 `
    1.93 "14/09/2023 23:13" 2.4 "21/09/2023 02:22" `
    3.55 "14/09/2023 23:13" 3.22 "21/09/2023 02:22" `
    4.3 "14/09/2023 23:13" 3.3 "21/09/2023 02:28" `
    "https://www.betexplorer.com/football/brazil/serie-a/sao-paulo-fortaleza/fHxE3NGR/"

Write-Host "Edit complete"
